$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 2 (shifts existing rows 2-18 down to 3-19)
$ws.Rows.Item(2).Insert()

# Copy the date-format style from the old row 2's A cell (now A3) into new A2 only
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Undo any style bleed onto B2:E2 caused by the row insert/paste - clear to default formatting
$ws.Range("B2:E2").ClearFormats()

# Full target data for rows 2..19: Date, y_0, y_0_forecast, y_1, y_1_forecast
$data = @(
    @(39400, 2007, 2.070003986395053, 2008, 0.6967455006573253),
    @(39765, 2008, 0.517569958955022, 2009, -0.6367039903686034),
    @(40130, 2009, -3.956152295564896, 2010, -0.6950853968889392),
    @(40494, 2010, 1.234995474941392, 2011, -0.2098161877568061),
    @(40862, 2011, 0.899360810820804, 2012, 1.205741443109987),
    @(41228, 2012, 0.9010266119894084, 2013, 1.531699207045123),
    @(41592, 2013, 0.02019328874804938, 2014, -0.1259279434590921),
    @(41957, 2014, 0.1729981757035093, 2015, 0.1749537368921361),
    @(42321, 2015, 0.09752710595589686, 2016, -0.001769149545471915),
    @(42689, 2016, -0.5280591151586633, 2017, -0.05116199209030947),
    @(43053, 2017, 0.07201851318385799, 2018, 0.2735900898381383),
    @(43418, 2018, 0.3727661260635617, 2019, -0.9505847809128332),
    @(43783, 2019, -0.801759526476209, 2020, 0.047674034857903),
    @(44159, 2020, -1.103489789942047, 2021, 1.605918384453009),
    @(44525, 2021, 0.9704846793491928, 2022, -0.8255212498362474),
    @(44890, 2022, -0.7009264669202708, 2023, 0.6624163082313173),
    @(45254, 2023, 0.3928252664241905, 2024, 0.302295480375836),
    @(45618, 2024, 0.3224026462283813, 2025, -0.7618983399156787)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
